$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels
$ws.Range("A1").Value = "peak_age"
$ws.Range("B1").Value = "peak_wage"
$ws.Range("C1").Value = "ci_wage_lo"
$ws.Range("D1").Value = "ci_wage_hi"

# Apply the same bold/centered header formatting used by the existing headers
$ws.Range("C1:D1").Font.Bold = $true
$ws.Range("C1:D1").HorizontalAlignment = -4108

# Update data row values
$ws.Range("B2").Value = 6259.493667603875
$ws.Range("C2").Value = 6133.368678375889
$ws.Range("D2").Value = 6389.226704294119
